# Update Excel SCD0011 until SCD0016
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0176 to SCD0011
$ws.Name = "SCD0011"

# Update cell B2 value (TC_ID) from DGS-191 to SCD0011-007
$ws.Range("B2").Value = "SCD0011-007"

# Column B width changed (widened to fit the longer TC id)
$ws.Columns("B").ColumnWidth = 11.67

# Selection / view changes: scroll back to show column A (remove topLeftCell=K1),
# and move active selection from L3 to B3
$ws.Range("B3").Select() | Out-Null
